# Insert a new weekly record at row 73 (pushing all subsequent rows down by one)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(73).Insert()

$ws.Range("A73").Value = 3
$ws.Range("B73").Value = "Femacal de La Calera"
$ws.Range("C73").Value = "Coquimbo"
$ws.Range("D73").Value = 44452
$ws.Range("E73").Value = 5
$ws.Range("F73").Value = 100112010
$ws.Range("G73").Value = "Achicoria"
$ws.Range("H73").Value = "Sin especificar"
$ws.Range("I73").Value = "Primera"
$ws.Range("J73").Value = 130
$ws.Range("K73").Value = 6500
$ws.Range("L73").Value = 7000
$ws.Range("M73").Value = 6731
$ws.Range("N73").Value = "$/caja 16 unidades"
$ws.Range("O73").Value = "Provincia de Quillota"
$ws.Range("P73").Value = 421
$ws.Range("Q73").Value = 16
$ws.Range("R73").Value = "Hortaliza"
